$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H29").Value = 1299.6666
$ws.Range("J29").Value = 3750
$ws.Range("L29").Value = 11250
$ws.Range("N29").Value = -11812
$ws.Range("H38").Value = 2478.875
$ws.Range("J38").Value = 6336
$ws.Range("L38").Value = 19008
$ws.Range("N38").Value = -19752
$ws.Range("H58").Value = 5374.5
$ws.Range("I58").Value = 750
$ws.Range("J58").Value = 9999
$ws.Range("K58").Value = 2250
$ws.Range("L58").Value = 29997
$ws.Range("M58").Value = -2100
$ws.Range("N58").Value = -30297
$ws.Range("H87").Value = 92677
$ws.Range("J87").Value = 95354
$ws.Range("L87").Value = 95354
$ws.Range("N87").Value = -97850
$ws.Range("H90").Value = 92677
$ws.Range("J90").Value = 95354
$ws.Range("L90").Value = 286062
$ws.Range("N90").Value = -298542
$ws.Range("H97").Value = 449.66666
$ws.Range("J97").Value = 449.66666
$ws.Range("L97").Value = 1348.99998
$ws.Range("N97").Value = -2340.99998
$ws.Range("H112").Value = 1315.0714
$ws.Range("J112").Value = 1657
$ws.Range("L112").Value = 4971
$ws.Range("N112").Value = -7187
$ws.Range("H116").Value = 4857.222
$ws.Range("I116").Value = 4999
$ws.Range("J116").Value = 4573.6665
$ws.Range("K116").Value = 4999
$ws.Range("L116").Value = 4573.6665
$ws.Range("M116").Value = -1557
$ws.Range("N116").Value = -11457.6665
$ws.Range("H125").Value = 2040.4286
$ws.Range("I125").Value = 1932.5
$ws.Range("J125").Value = 2688
$ws.Range("K125").Value = 17392.5
$ws.Range("L125").Value = 24192
$ws.Range("M125").Value = -14932.5
$ws.Range("N125").Value = -29112
$ws.Range("H135").Value = 750.2105
$ws.Range("I135").Value = 719.625
$ws.Range("J135").Value = 913.3333
$ws.Range("K135").Value = 6476.625
$ws.Range("L135").Value = 8219.9997
$ws.Range("M135").Value = -3941.625
$ws.Range("N135").Value = -13289.9997
$ws.Range("H137").Value = 3801.818
$ws.Range("I137").Value = 1984.4
$ws.Range("K137").Value = 5953.200000000001
$ws.Range("M137").Value = -3403.200000000001
$ws.Range("H138").Value = 7634.0625
$ws.Range("I138").Value = 1349.2858
$ws.Range("K138").Value = 4047.8574
$ws.Range("M138").Value = 1092.1426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7351.2666
$ws.Range("I32").Value = 7273.3335
$ws.Range("K32").Value = 7273.3335
$ws.Range("M32").Value = -6986.3335
$ws.Range("H45").Value = 1262
$ws.Range("I45").Value = 1262
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1262
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -885
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 2996.5
$ws.Range("I61").Value = 2996
$ws.Range("K61").Value = 2996
$ws.Range("M61").Value = -2784
$ws.Range("H74").Value = 1282.6666
$ws.Range("I74").Value = 1500
$ws.Range("K74").Value = 1500
$ws.Range("M74").Value = -626
$ws.Range("H77").Value = 1282.6666
$ws.Range("I77").Value = 1500
$ws.Range("K77").Value = 7500
$ws.Range("M77").Value = -3132
$ws.Range("H132").Value = 2381.56
$ws.Range("I132").Value = 1932.9445
$ws.Range("K132").Value = 5798.833500000001
$ws.Range("M132").Value = -3268.833500000001
$ws.Range("H136").Value = 2996.5
$ws.Range("I136").Value = 2996
$ws.Range("K136").Value = 8988
$ws.Range("M136").Value = -6438

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 595
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H64").Value = 1999.5
$ws.Range("J64").Value = 1999.5
$ws.Range("L64").Value = 1999.5
$ws.Range("N64").Value = -2449.5
$ws.Range("H67").Value = 1999.5
$ws.Range("J67").Value = 1999.5
$ws.Range("L67").Value = 1999.5
$ws.Range("N67").Value = -3559.5
$ws.Range("H86").Value = 407
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 407
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1693.9412
$ws.Range("I31").Value = 1791.6666
$ws.Range("K31").Value = 1791.6666
$ws.Range("M31").Value = -1496.6666
$ws.Range("H34").Value = 1693.9412
$ws.Range("I34").Value = 1791.6666
$ws.Range("K34").Value = 1791.6666
$ws.Range("M34").Value = -1589.6666
$ws.Range("H105").Value = 3474.5
$ws.Range("I105").Value = 2966
$ws.Range("K105").Value = 2966
$ws.Range("M105").Value = -1219
$ws.Range("H122").Value = 2916.2222
$ws.Range("I122").Value = 3093.25
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 9279.75
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -6829.75
$ws.Range("N122").Value = -9400
$ws.Range("H134").Value = 4276.8945
$ws.Range("I134").Value = 4309.5293
$ws.Range("K134").Value = 12928.5879
$ws.Range("M134").Value = -10393.5879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 3000
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 3000
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7169
$ws.Range("I70").Value = 7169
$ws.Range("K70").Value = 7169
$ws.Range("M70").Value = -6899
$ws.Range("H73").Value = 7169
$ws.Range("I73").Value = 7169
$ws.Range("K73").Value = 7169
$ws.Range("M73").Value = -6233
$ws.Range("H132").Value = 3320.5557
$ws.Range("I132").Value = 2647.8333
$ws.Range("K132").Value = 7943.499899999999
$ws.Range("M132").Value = -5413.499899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1051
$ws.Range("I7").Value = 901.3333
$ws.Range("J7").Value = 1500
$ws.Range("K7").Value = 901.3333
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = -789.3333
$ws.Range("N7").Value = -1724
$ws.Range("H40").Value = 5084.769
$ws.Range("I40").Value = 5439.091
$ws.Range("K40").Value = 5439.091
$ws.Range("M40").Value = -5303.091
$ws.Range("H53").Value = 10046
$ws.Range("I53").Value = 10046
$ws.Range("K53").Value = 10046
$ws.Range("M53").Value = -9528
$ws.Range("H55").Value = 283.42856
$ws.Range("I55").Value = 300
$ws.Range("J55").Value = 280.66666
$ws.Range("K55").Value = 300
$ws.Range("L55").Value = 280.66666
$ws.Range("M55").Value = -127
$ws.Range("N55").Value = -626.66666
$ws.Range("H126").Value = 1051
$ws.Range("I126").Value = 901.3333
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 2703.9999
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -233.9998999999998
$ws.Range("N126").Value = -9440
$ws.Range("H136").Value = 3654
$ws.Range("I136").Value = 3654
$ws.Range("K136").Value = 10962
$ws.Range("M136").Value = -8412

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3098
$ws.Range("I126").Value = 2830
$ws.Range("K126").Value = 8490
$ws.Range("M126").Value = -6020
$ws.Range("H132").Value = 2623.875
$ws.Range("I132").Value = 1478.0714
$ws.Range("J132").Value = 4228
$ws.Range("K132").Value = 4434.2142
$ws.Range("L132").Value = 12684
$ws.Range("M132").Value = -1904.2142
$ws.Range("N132").Value = -17744
$ws.Range("H136").Value = 628.2105
$ws.Range("I136").Value = 615.125
$ws.Range("K136").Value = 1845.375
$ws.Range("M136").Value = 704.625
